$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the RF (raising factor) column I for rows 28 through 38
# from the old value (100.1797142857143) to the new value (22.83225)
$ws.Range("I28:I38").Value = 22.83225
